# Updated cryptos list on Sun Aug 11 08:29:56 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a Price-column value while forcing text storage so that
# purely-numeric-looking strings (e.g. "6.20", "1.50") keep their exact
# textual representation instead of being auto-coerced into a Number by
# Excel's smart entry (mirrors the source file's inlineStr text cells).
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.176.20"
$ws.Range("E2").Value = "  +0.97%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.662.47"
$ws.Range("E3").Value = "  +1.82%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "533.45"
$ws.Range("E5").Value = "  +4.25%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "156.93"
$ws.Range("E6").Value = "  +1.39%  "

# Row 7 - USDC
Set-TextValue $ws.Range("D7") "0.997"
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.88%  "

# Row 9 - Toncoin
Set-TextValue $ws.Range("D9") "6.62"
$ws.Range("E9").Value = "  -1.91%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +5.08%  "

# Row 11 - Cardano
Set-TextValue $ws.Range("D11") "0.354"
$ws.Range("E11").Value = "  +2.22%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.127.80"
$ws.Range("E13").Value = "  +1.82%  "

# Row 14 - WrappedBTC
$ws.Range("D14").Value = "61.170.62"
$ws.Range("E14").Value = "  +1.08%  "

# Row 15 - Avalanche
Set-TextValue $ws.Range("D15") "22.17"
$ws.Range("E15").Value = "  +2.60%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +2.51%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.670.56"
$ws.Range("E17").Value = "  +1.87%  "

# Row 18 - Polkadot
Set-TextValue $ws.Range("D18") "4.79"
$ws.Range("E18").Value = "  +0.70%  "

# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "355.81"
$ws.Range("E19").Value = "  +0.92%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "10.74"
$ws.Range("E20").Value = "  +1.36%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "6.31"
$ws.Range("E21").Value = "  +2.12%  "

# Row 22 - Dai
Set-TextValue $ws.Range("D22") "0.999"
$ws.Range("E22").Value = "  -0.01%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "61.67"
$ws.Range("E23").Value = "  +1.69%  "

# Row 24 - Polygon
Set-TextValue $ws.Range("D24") "0.434"
$ws.Range("E24").Value = "  +2.46%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  +1.54%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  +0.68%  "

# Row 27 - PEPE (contains subscript-3 U+2083)
$sub3 = [char]0x2083
$ws.Range("D27").Value = "0.0$($sub3)0865"
$ws.Range("E27").Value = "  +2.73%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D28") "7.43"
$ws.Range("E28").Value = "  +1.12%  "

# Row 29 - USDe
$ws.Range("E29").Value = "  -0.04%  "

# Row 30 - Aptos
Set-TextValue $ws.Range("D30") "6.20"
$ws.Range("E30").Value = "  +6.58%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +4.02%  "

# Row 32 - EthereumClassic
Set-TextValue $ws.Range("D32") "19.61"
$ws.Range("E32").Value = "  +0.90%  "

# Row 33 - Monero
Set-TextValue $ws.Range("D33") "150.60"
$ws.Range("E33").Value = "  -0.52%  "

# Row 34 - NEARProtocol
Set-TextValue $ws.Range("D34") "4.17"
$ws.Range("E34").Value = "  +4.62%  "

# Row 35 - ImmutableX
Set-TextValue $ws.Range("D35") "1.21"
$ws.Range("E35").Value = "  +1.37%  "

# Row 36 - Fetch.AI
Set-TextValue $ws.Range("D36") "0.922"
$ws.Range("E36").Value = "  +9.05%  "

# Row 37 - SuiNetwork
Set-TextValue $ws.Range("D37") "0.885"
$ws.Range("E37").Value = "  -0.79%  "

# Row 38 - was Stacks, now Bittensor
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D38") "309.14"
$ws.Range("E38").Value = "  +6.14%  "

# Row 39 - was Bittensor, now Stacks
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D39") "1.50"
$ws.Range("E39").Value = "  +1.17%  "

# Row 40 - Filecoin
Set-TextValue $ws.Range("D40") "3.83"
$ws.Range("E40").Value = "  +1.98%  "

# Row 41 - Mantle
$ws.Range("E41").Value = "  +4.24%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  +2.10%  "

# Row 43 - Hedera
Set-TextValue $ws.Range("D43") "0.0568"
$ws.Range("E43").Value = "  +2.41%  "

# Row 44 - EnergySwap
Set-TextValue $ws.Range("D44") "20.34"
$ws.Range("E44").Value = "  +2.37%  "

# Row 46 - RenderToken
Set-TextValue $ws.Range("D46") "5.05"
$ws.Range("E46").Value = "  +2.76%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  +2.40%  "

# Row 48 - InjectiveProtocol
$ws.Range("E48").Value = "  +9.12%  "

# Row 49 - WhiteBITCoin
Set-TextValue $ws.Range("D49") "10.37"
$ws.Range("E49").Value = "  +0.50%  "

# Row 50 - Maker
$ws.Range("D50").Value = "2.002.14"
$ws.Range("E50").Value = "  +0.15%  "

# Row 51 - dogwifhat
Set-TextValue $ws.Range("D51") "1.86"
$ws.Range("E51").Value = "  +2.85%  "
